$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.124.21'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.789.49'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.546'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.36'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.296'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0688'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0939'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('D12').Value = '2.047.28'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.52'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.92%  '
$ws.Range('D14').Value = '1.775.46'
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').Value = '34.103.18'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').Value = '0.0₃0778'
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.05%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.04'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.12%  '
$ws.Range('E26').Value = '  +1.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.26'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  +0.90%  '
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0518'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('E33').Value = '  +3.22%  '
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('D35').Value = '1.403.94'
$ws.Range('E35').Value = '  +0.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.651'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('E37').Value = '  +1.81%  '
$ws.Range('E38').Value = '  -0.83%  '
$ws.Range('E39').Value = '  +6.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.38'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '80.08'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.922'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.50%  '
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.38'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +10.96%  '
$ws.Range('D45').Value = '0.0₆0140'
$ws.Range('E45').Value = '  -0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.81%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.08'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.40%  '
$ws.Range('B48').Value = 'Kaspa'
$ws.Range('C48').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0507'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '106.97'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('D50').Value = '1.948.87'
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('E51').Value = '  +0.15%  '
